# Added 4wk low sales check
# Updates forecast (MyForecast), Inventory Coverage, and Seasonality Index
# values on the "Forecast Comparison" sheet, and refreshes the summary
# metrics on the "Summary" sheet to match.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# Row -> (MyForecast, Inventory Coverage, Seasonality Index)
$rows = @(
    @{ Row = 2;  D = 2; H = 32;                  L = 1.18 },
    @{ Row = 3;  D = 3; H = 23.48;                L = 0.9 },
    @{ Row = 4;  D = 3; H = 20.05;                L = 0.8100000000000001 },
    @{ Row = 5;  D = 3; H = 21.36;                L = 1.01 },
    @{ Row = 6;  D = 2; H = 29.22;                L = 1.04 },
    @{ Row = 7;  D = 1; H = 43.27;                L = 0.84 },
    @{ Row = 8;  D = 1; H = 57.64;                L = 0.89 },
    @{ Row = 9;  D = 0; H = 69.22;                L = 0.9399999999999999 },
    @{ Row = 10; D = 0; H = 68.22;                L = 0.85 },
    @{ Row = 11; D = 0; H = 86.43000000000001;    L = 1 },
    @{ Row = 12; D = 0; H = 119.6;                L = 0.83 },
    @{ Row = 13; D = 0; H = 118.6;                L = 0.93 },
    @{ Row = 14; D = 0; H = 117.6;                L = 1.08 },
    @{ Row = 15; D = 0; H = 83.29000000000001;    L = 0.83 },
    @{ Row = 16; D = 0; H = 115.2;                L = 0.93 },
    @{ Row = 17; D = 0; H = 114.2;                L = 1.19 }
)

foreach ($r in $rows) {
    $wsForecast.Cells.Item($r.Row, 4).Value = $r.D
    $wsForecast.Cells.Item($r.Row, 8).Value = $r.H
    $wsForecast.Cells.Item($r.Row, 12).Value = $r.L
}

$summaryUpdates = @(
    @{ Cell = "B9";  Text = "23" },
    @{ Cell = "B10"; Text = "19" },
    @{ Cell = "B11"; Text = "13" },
    @{ Cell = "B12"; Text = "4" },
    @{ Cell = "B14"; Text = "0" }
)

foreach ($u in $summaryUpdates) {
    $cell = $wsSummary.Range($u.Cell)
    # Force the numeric-looking string to be stored as text (matching the
    # original inline-string cell type) without leaving a residual
    # non-default cell style behind.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Text
    $cell.Style = "Normal"
}
